$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix naive component forecaster bug - update QoQ error values for rows 25-53

# Row 25
$ws.Cells.Item(25, 10).Value = -0.9790448517306776
$ws.Cells.Item(25, 11).Value = 0.3108804844682282

# Row 26
$ws.Cells.Item(26, 9).Value = -1.279044851730669
$ws.Cells.Item(26, 10).Value = 0.01088048446822254

# Row 27
$ws.Cells.Item(27, 8).Value = -1.329044851730686
$ws.Cells.Item(27, 9).Value = -0.03911951553179449

# Row 28
$ws.Cells.Item(28, 7).Value = -1.479044851730669
$ws.Cells.Item(28, 8).Value = -0.1891195155317774

# Row 29
$ws.Cells.Item(29, 6).Value = -1.179044851730652
$ws.Cells.Item(29, 7).Value = 0.1108804844681828
$ws.Cells.Item(29, 8).Value = 2.32007603651995
$ws.Cells.Item(29, 9).Value = -0.585713798234863
$ws.Cells.Item(29, 10).Value = -0.2260365910423587
$ws.Cells.Item(29, 11).Value = -0.6238769106320059

# Row 30
$ws.Cells.Item(30, 5).Value = -1.228044851730673
$ws.Cells.Item(30, 6).Value = 0.06188048446821881
$ws.Cells.Item(30, 7).Value = 2.271076036519957
$ws.Cells.Item(30, 8).Value = -0.634713798234827
$ws.Cells.Item(30, 9).Value = -0.2750365910423511
$ws.Cells.Item(30, 10).Value = -0.6728769106320267

# Row 31
$ws.Cells.Item(31, 4).Value = -1.17904485173068
$ws.Cells.Item(31, 5).Value = 0.1108804844682396
$ws.Cells.Item(31, 6).Value = 2.320076036519907
$ws.Cells.Item(31, 7).Value = -0.585713798234792
$ws.Cells.Item(31, 8).Value = -0.2260365910423872
$ws.Cells.Item(31, 9).Value = -0.6238769106320627

# Row 32
$ws.Cells.Item(32, 3).Value = -0.879044851730669
$ws.Cells.Item(32, 4).Value = 0.2108804844682226
$ws.Cells.Item(32, 5).Value = 2.220076036519961
$ws.Cells.Item(32, 6).Value = -0.685713798234809
$ws.Cells.Item(32, 7).Value = -0.3260365910423758
$ws.Cells.Item(32, 8).Value = -0.723876910632023

# Row 33
$ws.Cells.Item(33, 2).Value = -1.080860442277157
$ws.Cells.Item(33, 3).Value = 0.2108804844682055
$ws.Cells.Item(33, 4).Value = 2.320344273492722
$ws.Cells.Item(33, 5).Value = -0.5857137982348204
$ws.Cells.Item(33, 6).Value = -0.2260365910423872
$ws.Cells.Item(33, 7).Value = -0.6238769106320627
$ws.Cells.Item(33, 8).Value = 2.836777872354247
$ws.Cells.Item(33, 9).Value = -0.09340513149564345
$ws.Cells.Item(33, 10).Value = -0.3299754819599623
$ws.Cells.Item(33, 11).Value = -0.004381049978860574

# Row 34
$ws.Cells.Item(34, 2).Value = 0.2080883825945676
$ws.Cells.Item(34, 3).Value = 2.320384586168557
$ws.Cells.Item(34, 4).Value = -0.5857137982348204
$ws.Cells.Item(34, 5).Value = -0.2260365910424156
$ws.Cells.Item(34, 6).Value = -0.6238769106320343
$ws.Cells.Item(34, 7).Value = 2.836777872354205
$ws.Cells.Item(34, 8).Value = -0.09340513149561502
$ws.Cells.Item(34, 9).Value = -0.3299754819600049
$ws.Cells.Item(34, 10).Value = -0.004381049978888996

# Row 35
$ws.Cells.Item(35, 2).Value = 2.320739296142662
$ws.Cells.Item(35, 3).Value = -0.5857137982348204
$ws.Cells.Item(35, 4).Value = -0.2260365910424156
$ws.Cells.Item(35, 5).Value = -0.6238769106320059
$ws.Cells.Item(35, 6).Value = 2.836777872354233
$ws.Cells.Item(35, 7).Value = -0.09340513149561502
$ws.Cells.Item(35, 8).Value = -0.3299754819599907
$ws.Cells.Item(35, 9).Value = -0.004381049978860574

# Row 36
$ws.Cells.Item(36, 2).Value = -0.5851877673051857
$ws.Cells.Item(36, 3).Value = -0.2260365910423872
$ws.Cells.Item(36, 4).Value = -0.6238769106320343
$ws.Cells.Item(36, 5).Value = 2.836777872354205
$ws.Cells.Item(36, 6).Value = -0.09340513149567187
$ws.Cells.Item(36, 7).Value = -0.3299754819599623
$ws.Cells.Item(36, 8).Value = -0.004381049978860574

# Row 37
$ws.Cells.Item(37, 2).Value = 0.5758294174215872
$ws.Cells.Item(37, 3).Value = 0.176123089368005
$ws.Cells.Item(37, 4).Value = 3.835785141722042
$ws.Cells.Item(37, 5).Value = 0.7065948685043819
$ws.Cells.Item(37, 6).Value = -0.02997548195997934
$ws.Cells.Item(37, 7).Value = 0.2956189500211224
$ws.Cells.Item(37, 8).Value = 0.7593444227005932
$ws.Cells.Item(37, 9).Value = 0.473188827696518
$ws.Cells.Item(37, 10).Value = 0.5950343932174889
$ws.Cells.Item(37, 11).Value = 0.1613548489792436

# Row 38
$ws.Cells.Item(38, 2).Value = 0.176834182002523
$ws.Cells.Item(38, 3).Value = 3.834372229001446
$ws.Cells.Item(38, 4).Value = 0.7065948685043959
$ws.Cells.Item(38, 5).Value = -0.02997548195999356
$ws.Cells.Item(38, 6).Value = 0.2956189500211654
$ws.Cells.Item(38, 7).Value = 0.7593444227005364
$ws.Cells.Item(38, 8).Value = 0.4731888276965322
$ws.Cells.Item(38, 9).Value = 0.5950343932175173
$ws.Cells.Item(38, 10).Value = 0.1613548489792436

# Row 39
$ws.Cells.Item(39, 2).Value = 3.534544165640355
$ws.Cells.Item(39, 3).Value = 0.906594868504385
$ws.Cells.Item(39, 4).Value = -0.02997548195997934
$ws.Cells.Item(39, 5).Value = 0.1956189500211281
$ws.Cells.Item(39, 6).Value = 0.7593444227005506
$ws.Cells.Item(39, 7).Value = 0.4731888276965322
$ws.Cells.Item(39, 8).Value = 0.1950343932175116
$ws.Cells.Item(39, 9).Value = -0.2386451510207763

# Row 40
$ws.Cells.Item(40, 2).Value = 0.1079662819227423
$ws.Cells.Item(40, 3).Value = 0.1700245180400941
$ws.Cells.Item(40, 4).Value = 0.4956189500210963
$ws.Cells.Item(40, 5).Value = 0.8593444227005591
$ws.Cells.Item(40, 6).Value = 0.4731888276965322
$ws.Cells.Item(40, 7).Value = 0.4950343932175088
$ws.Cells.Item(40, 8).Value = 0.0613548489793061

# Row 41
$ws.Cells.Item(41, 2).Value = -0.6261225043750751
$ws.Cells.Item(41, 3).Value = -0.3043810499788719
$ws.Cells.Item(41, 4).Value = 0.658702131206482
$ws.Cells.Item(41, 5).Value = 0.3732545146475133
$ws.Cells.Item(41, 6).Value = 0.6951662308824917
$ws.Cells.Item(41, 7).Value = 0.2615081725498529
$ws.Cells.Item(41, 8).Value = 2.147035121354733
$ws.Cells.Item(41, 9).Value = 0.2436720223698501
$ws.Cells.Item(41, 10).Value = -0.1119095091543443
$ws.Cells.Item(41, 11).Value = 0.675040353419746

# Row 42
$ws.Cells.Item(42, 2).Value = 0.294793270082792
$ws.Cells.Item(42, 3).Value = 0.7543711825399271
$ws.Cells.Item(42, 4).Value = 0.4735851177026121
$ws.Cells.Item(42, 5).Value = 0.5954391506335384
$ws.Cells.Item(42, 6).Value = 0.06189789600354345
$ws.Cells.Item(42, 7).Value = 1.946527446418656
$ws.Cells.Item(42, 8).Value = 0.2436720223699353
$ws.Cells.Item(42, 9).Value = 0.1880904908456245
$ws.Cells.Item(42, 10).Value = 0.9750403534197574

# Row 43
$ws.Cells.Item(43, 2).Value = 0.4943829294508504
$ws.Cells.Item(43, 3).Value = 0.2913897688964938
$ws.Cells.Item(43, 4).Value = 0.497312044872146
$ws.Cells.Item(43, 5).Value = 0.1780735536881224
$ws.Cells.Item(43, 6).Value = 1.980143888005614
$ws.Cells.Item(43, 7).Value = 0.2584201545070783
$ws.Cells.Item(43, 8).Value = 0.2030944081333246
$ws.Cells.Item(43, 9).Value = 0.9940013464764803

# Row 44
$ws.Cells.Item(44, 2).Value = -0.2133677787764014
$ws.Cells.Item(44, 3).Value = 0.1617784409566667
$ws.Cells.Item(44, 4).Value = 0.5963864610890259
$ws.Cells.Item(44, 5).Value = 2.313891085639355
$ws.Cells.Item(44, 6).Value = 0.426043548304591
$ws.Cells.Item(44, 7).Value = 0.2876270601879725
$ws.Cells.Item(44, 8).Value = 1.082704344158415

# Row 45
$ws.Cells.Item(45, 2).Value = 0.0579608153039004
$ws.Cells.Item(45, 3).Value = 0.2256477827406087
$ws.Cells.Item(45, 4).Value = 2.484896471875359
$ws.Cells.Item(45, 5).Value = 0.4437418058852245
$ws.Cells.Item(45, 6).Value = 0.3359833184482852
$ws.Cells.Item(45, 7).Value = 0.939729490273109
$ws.Cells.Item(45, 8).Value = -1.415398033202167
$ws.Cells.Item(45, 9).Value = -0.659337854358782

# Row 46
$ws.Cells.Item(46, 2).Value = -0.4295258376674695
$ws.Cells.Item(46, 3).Value = 1.774655533798551
$ws.Cells.Item(46, 4).Value = 0.04286727722244166
$ws.Cells.Item(46, 5).Value = -0.01237347804517697
$ws.Cells.Item(46, 6).Value = 0.775168870215893
$ws.Cells.Item(46, 7).Value = -1.414323171601396
$ws.Cells.Item(46, 8).Value = -0.5766660043058638

# Row 47
$ws.Cells.Item(47, 2).Value = 1.781099436349905
$ws.Cells.Item(47, 3).Value = 0.0450760692686174
$ws.Cells.Item(47, 4).Value = -0.0111680561212637
$ws.Cells.Item(47, 5).Value = 0.779269356066564
$ws.Cells.Item(47, 6).Value = -1.418535507728536
$ws.Cells.Item(47, 7).Value = -0.5766660043059064

# Row 48
$ws.Cells.Item(48, 2).Value = -0.2578266341250811
$ws.Cells.Item(48, 3).Value = -0.2119095091543527
$ws.Cells.Item(48, 4).Value = 0.9403763896489465
$ws.Cells.Item(48, 5).Value = -1.251139686906342
$ws.Cells.Item(48, 6).Value = -0.4944766814167993

# Row 49
$ws.Cells.Item(49, 2).Value = -0.4069492870295619
$ws.Cells.Item(49, 3).Value = 0.6750403534197744
$ws.Cells.Item(49, 4).Value = -1.220764540768429
$ws.Cells.Item(49, 5).Value = -0.3766660043058609

# Row 50
$ws.Cells.Item(50, 2).Value = 0.8628474748582136
$ws.Cells.Item(50, 3).Value = -1.339093115628018
$ws.Cells.Item(50, 4).Value = -0.5920667263155368

# Row 51
$ws.Cells.Item(51, 2).Value = -1.436438518536832
$ws.Cells.Item(51, 3).Value = -0.5918833611148815

# Row 52
$ws.Cells.Item(52, 2).Value = -0.6276069079710285

# Clear trailing cells that are no longer populated
$ws.Cells.Item(45, 10).ClearContents()
$ws.Cells.Item(46, 9).ClearContents()
$ws.Cells.Item(47, 8).ClearContents()
$ws.Cells.Item(48, 7).ClearContents()
$ws.Cells.Item(49, 6).ClearContents()
$ws.Cells.Item(50, 5).ClearContents()
$ws.Cells.Item(51, 4).ClearContents()
$ws.Cells.Item(52, 3).ClearContents()
$ws.Cells.Item(53, 2).ClearContents()

